$d = $word.ActiveDocument

# 1) Remove the two comments (id=0 "que hidden..." and id=1 "uille hidden...").
#    Deleting each Comment removes its commentRangeStart/commentRangeEnd/
#    commentReference markup around "q" and "u" in the body text, leaving the
#    plain runs "q" and "u" behind untouched.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# 2) Fix small transcription typos in the body text.
$d.Content.Find.Execute("contenir aveclq du", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "contenir avecq du", 2) | Out-Null

$d.Content.Find.Execute("fil et", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "filet", 2) | Out-Null

$d.Save()
